$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the bank journal names to include IBAN hints
$ws.Range("C2").Value = "B. Pop. Software (IT15*456)"
$ws.Range("C3").Value = "B. Credito per Tutti (IT74*680)"
$ws.Range("C4").Value = "Portafoglio RiBA B. Pop. Soft. (IT15*456/IT26*456)"
$ws.Range("C5").Value = "Portafoglio Anticipi B. Pop. Soft. (IT15*456/IT82*456)"
$ws.Range("C6").Value = "Portafoglio Anticipi Estero BCT (IT74*680/IT58*680)"

# Widen column C to fit the new longer text (target ~46.12 characters;
# COM ColumnWidth snaps to the nearest whole pixel, so 45.25 is the input
# that lands closest to the desired stored width)
$ws.Columns("C").ColumnWidth = 45.25

# Move the active selection to C6
$ws.Range("C6").Select()
